# Add the new "2025-07" monthly data row (row 20) to each of the
# 11 city/prefecture worksheets that still only have rows 1-19.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("杭州市")
$ws.Range("A20").Value = "2025-07"
$ws.Range("B20").Value = 7877947.975100003
$ws.Range("C20").Value = 0.0751204891831696
$ws.Range("D20").Value = 51537900.0045
$ws.Range("E20").Value = 0.07461906266714435
$ws.Range("F20").Value = 5813005.852599997
$ws.Range("G20").Value = 0.1146547304019987
$ws.Range("H20").Value = 36795478.9278
$ws.Range("I20").Value = 0.1284437213123026
$ws.Range("J20").Value = 2064942.122499999
$ws.Range("K20").Value = -0.0224799138707521
$ws.Range("L20").Value = 14742421.0767
$ws.Range("M20").Value = -0.03970359956219227

$ws = $wb.Worksheets.Item("湖州市")
$ws.Range("A20").Value = "2025-07"
$ws.Range("B20").Value = 1288651.6972
$ws.Range("C20").Value = 0.07292555650091104
$ws.Range("D20").Value = 10502381.4693
$ws.Range("E20").Value = 0.04479370048171738
$ws.Range("F20").Value = 1112730.2829
$ws.Range("G20").Value = 0.03976943021706325
$ws.Range("H20").Value = 9728602.788700001
$ws.Range("I20").Value = 0.04946209765809084
$ws.Range("J20").Value = 175921.4143000001
$ws.Range("K20").Value = 0.3440067506306597
$ws.Range("L20").Value = 773778.6806000001
$ws.Range("M20").Value = -0.0105452223601773

$ws = $wb.Worksheets.Item("嘉兴市")
$ws.Range("A20").Value = "2025-07"
$ws.Range("B20").Value = 3939099.797099996
$ws.Range("C20").Value = -0.0285426966504938
$ws.Range("D20").Value = 28350864.216
$ws.Range("E20").Value = 0.03813571396607318
$ws.Range("F20").Value = 2914879.941799998
$ws.Range("G20").Value = -0.06475928277018894
$ws.Range("H20").Value = 21618939.7088
$ws.Range("I20").Value = 0.03514497942502848
$ws.Range("J20").Value = 1024219.8553
$ws.Range("K20").Value = 0.09177973920341298
$ws.Range("L20").Value = 6731924.5072
$ws.Range("M20").Value = 0.04785813284547946

$ws = $wb.Worksheets.Item("金华市")
$ws.Range("A20").Value = "2025-07"
$ws.Range("B20").Value = 8942821.705200002
$ws.Range("C20").Value = 0.2162906205160331
$ws.Range("D20").Value = 59810684.07080001
$ws.Range("E20").Value = 0.2031415980955018
$ws.Range("F20").Value = 7774597.428199999
$ws.Range("G20").Value = 0.2083787868046469
$ws.Range("H20").Value = 53111170.4147
$ws.Range("I20").Value = 0.206114335098879
$ws.Range("J20").Value = 1168224.277000001
$ws.Range("K20").Value = 0.27170361866487
$ws.Range("L20").Value = 6699513.6561
$ws.Range("M20").Value = 0.1800835072464051

$ws = $wb.Worksheets.Item("丽水市")
$ws.Range("A20").Value = "2025-07"
$ws.Range("B20").Value = 356033.8417000002
$ws.Range("C20").Value = -0.03984746131165684
$ws.Range("D20").Value = 2698698.4058
$ws.Range("E20").Value = 0.2283790325795219
$ws.Range("F20").Value = 306749.6035
$ws.Range("G20").Value = -0.09258057873725389
$ws.Range("H20").Value = 2339105.0071
$ws.Range("I20").Value = 0.2002954602207057
$ws.Range("J20").Value = 49284.23820000002
$ws.Range("K20").Value = 0.5042395166993727
$ws.Range("L20").Value = 359593.3987
$ws.Range("M20").Value = 0.4488944415765801

$ws = $wb.Worksheets.Item("宁波市")
$ws.Range("A20").Value = "2025-07"
$ws.Range("B20").Value = 12120659.72139999
$ws.Range("C20").Value = -0.03578740973043115
$ws.Range("D20").Value = 84300357.413
$ws.Range("E20").Value = 0.04480826275655891
$ws.Range("F20").Value = 8147799.170100003
$ws.Range("G20").Value = -0.01470602999827564
$ws.Range("H20").Value = 57191361.32380001
$ws.Range("I20").Value = 0.08315044783896575
$ws.Range("J20").Value = 3972860.5513
$ws.Range("K20").Value = -0.07631887189846964
$ws.Range("L20").Value = 27108996.0892
$ws.Range("M20").Value = -0.02779605578946598

$ws = $wb.Worksheets.Item("衢州市")
$ws.Range("A20").Value = "2025-07"
$ws.Range("B20").Value = 655915.2377999998
$ws.Range("C20").Value = -0.03621486849138078
$ws.Range("D20").Value = 4731846.2845
$ws.Range("E20").Value = 0.06596266197234213
$ws.Range("F20").Value = 468122.4135999996
$ws.Range("G20").Value = 0.02599909731569072
$ws.Range("H20").Value = 3388715.4907
$ws.Range("I20").Value = 0.1089227078320829
$ws.Range("J20").Value = 187792.8242000004
$ws.Range("K20").Value = -0.1627665502688909
$ws.Range("L20").Value = 1343130.7938
$ws.Range("M20").Value = -0.0289495420610969

$ws = $wb.Worksheets.Item("绍兴市")
$ws.Range("A20").Value = "2025-07"
$ws.Range("B20").Value = 3049188.866599999
$ws.Range("C20").Value = -0.005760055849616363
$ws.Range("D20").Value = 21620606.0246
$ws.Range("E20").Value = -0.1601843707521476
$ws.Range("F20").Value = 2880356.789900001
$ws.Range("G20").Value = 0.0188179652820617
$ws.Range("H20").Value = 20122432.0634
$ws.Range("I20").Value = -0.1444643877253412
$ws.Range("J20").Value = 168832.0766999999
$ws.Range("K20").Value = -0.2956485869742815
$ws.Range("L20").Value = 1498173.9612
$ws.Range("M20").Value = -0.326419164282361

$ws = $wb.Worksheets.Item("台州市")
$ws.Range("A20").Value = "2025-07"
$ws.Range("B20").Value = 2352457.749400001
$ws.Range("C20").Value = 0.02703428045880862
$ws.Range("D20").Value = 17129577.5234
$ws.Range("E20").Value = 0.1104688773715383
$ws.Range("F20").Value = 2124113.667199997
$ws.Range("G20").Value = 0.02523112266823668
$ws.Range("H20").Value = 15668645.1974
$ws.Range("I20").Value = 0.1147936285052429
$ws.Range("J20").Value = 228344.0821999998
$ws.Range("K20").Value = 0.04411668407799407
$ws.Range("L20").Value = 1460932.326
$ws.Range("M20").Value = 0.06611103340160618

$ws = $wb.Worksheets.Item("温州市")
$ws.Range("A20").Value = "2025-07"
$ws.Range("B20").Value = 2704837.497199999
$ws.Range("C20").Value = 0.02388217869541043
$ws.Range("D20").Value = 17660131.4852
$ws.Range("E20").Value = 0.0657141951665472
$ws.Range("F20").Value = 2409257.6712
$ws.Range("G20").Value = 0.08395975693881397
$ws.Range("H20").Value = 15283397.8529
$ws.Range("I20").Value = 0.09978540944268288
$ws.Range("J20").Value = 295579.8260000001
$ws.Range("K20").Value = -0.2947303833090155
$ws.Range("L20").Value = 2376733.6323
$ws.Range("M20").Value = -0.1113224849878335

$ws = $wb.Worksheets.Item("舟山市")
$ws.Range("A20").Value = "2025-07"
$ws.Range("B20").Value = 2786340.537699997
$ws.Range("C20").Value = -0.193792536397168
$ws.Range("D20").Value = 20999932.2359
$ws.Range("E20").Value = -0.01732457695855272
$ws.Range("F20").Value = 895317.3621000005
$ws.Range("G20").Value = -0.147196669930003
$ws.Range("H20").Value = 6941496.3062
$ws.Range("I20").Value = -0.03444850035550673
$ws.Range("J20").Value = 1891023.1756
$ws.Range("K20").Value = -0.2141223591524501
$ws.Range("L20").Value = 14058435.9297
$ws.Range("M20").Value = -0.0086434964357609
